$d = $word.ActiveDocument

# 1. Title: "Stipulazione Use Case" -> "Stipulazione dei Use Case" (with bookmark between "dei " and "Use Case")
$d.Content.Find.Execute("Stipulazione Use Case", $true, $false, $false, $false, $false, $true, 1, $false, "Stipulazione dei Use Case", 2)
